# Update the LR-pairs sheet with the new TPM-derived values.
#
# The "Sending cluster" (column A) for each 4-row block shifts by one
# position (ECs -> FAPs -> MuSCs -> Resolving-Mac, wrapping around), and
# all of the ligand/receptor/edge-derived metrics (columns E-J, M-T) are
# recomputed for the new TPM values. Columns B (Ligand symbol), C
# (Receptor symbol), D (Target cluster), K and L are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
  @{ row=2;  A="FAPs";          E=2; F=0.6666666666666666; G=0.6742976666666666; H=2.022893; I=0.3960292783506769;  J=0.3960292783506769;  M=75.93400333333334;  N=227.80201;  O=0.140568299724637;  P=0.140568299724637;  Q=51.20212126832555;  R=460.81909141493;   S=0.05566916229892963;  T=0.05566916229892963  }
  @{ row=3;  A="FAPs";          E=2; F=0.6666666666666666; G=0.6742976666666666; H=2.022893; I=0.3960292783506769;  J=0.3960292783506769;  M=202.4456936666667;  N=607.337081; O=0.3747655292237945;  P=0.3747655292237945;  Q=136.5086588661481;  R=1228.577929795333; S=0.1484181220892088;   T=0.1484181220892088   }
  @{ row=4;  A="FAPs";          E=2; F=0.6666666666666666; G=0.6742976666666666; H=2.022893; I=0.3960292783506769;  J=0.3960292783506769;  M=101.2433646666667;  N=303.730094; O=0.1874207470284575;  P=0.1874207470284576;  Q=68.26816456021577;  R=614.413481041942;  S=0.07422410319362481;  T=0.07422410319362482  }
  @{ row=5;  A="FAPs";          E=2; F=0.6666666666666666; G=0.6742976666666666; H=2.022893; I=0.3960292783506769;  J=0.3960292783506769;  M=160.56988;          N=481.70964;  O=0.297245424023111;   P=0.297245424023111;   Q=108.2718954209467;  R=974.4470587885199; S=0.1177178907689136;   T=0.1177178907689136   }
  @{ row=6;  A="MuSCs";         E=3; F=1;                  G=0.8625470000000001; H=2.587641; I=0.5065920925430184;  J=0.5065920925430184;  M=75.93400333333334;  N=227.80201;  O=0.140568299724637;  P=0.140568299724637;  Q=65.49664677315667;  R=589.46982095841;   S=0.07121078910271804;  T=0.07121078910271804  }
  @{ row=7;  A="MuSCs";         E=3; F=1;                  G=0.8625470000000001; H=2.587641; I=0.5065920925430184;  J=0.5065920925430184;  M=202.4456936666667;  N=607.337081; O=0.3747655292237945;  P=0.3747655292237945;  Q=174.6189257351024;  R=1571.570331615921; S=0.1898532536624737;   T=0.1898532536624737   }
  @{ row=8;  A="MuSCs";         E=3; F=1;                  G=0.8625470000000001; H=2.587641; I=0.5065920925430184;  J=0.5065920925430184;  M=101.2433646666667;  N=303.730094; O=0.1874207470284575;  P=0.1874207470284576;  Q=87.32716046313934;  R=785.944444168254;  S=0.094945868423122;    T=0.09494586842312201  }
  @{ row=9;  A="MuSCs";         E=3; F=1;                  G=0.8625470000000001; H=2.587641; I=0.5065920925430184;  J=0.5065920925430184;  M=160.56988;          N=481.70964;  O=0.297245424023111;   P=0.297245424023111;   Q=138.49906828436;    R=1246.49161455924;  S=0.1505821813547046;   T=0.1505821813547046   }
  @{ row=10; A="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.1658013333333333; H=0.497404; I=0.09737862910630474; J=0.09737862910630474; M=75.93400333333334;  N=227.80201;  O=0.140568299724637;  P=0.140568299724637;  Q=12.58995899800444;  R=113.30963098204;   S=0.0136883483229893;   T=0.0136883483229893   }
  @{ row=11; A="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.1658013333333333; H=0.497404; I=0.09737862910630474; J=0.09737862910630474; M=202.4456936666667;  N=607.337081; O=0.3747655292237945;  P=0.3747655292237945;  Q=33.56576593752489;  R=302.091893437724;  S=0.03649415347211189;  T=0.03649415347211189  }
  @{ row=12; A="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.1658013333333333; H=0.497404; I=0.09737862910630474; J=0.09737862910630474; M=101.2433646666667;  N=303.730094; O=0.1874207470284575;  P=0.1874207470284576;  Q=16.78628485288622;  R=151.076563675976;  S=0.01825077541171073;  T=0.01825077541171074  }
  @{ row=13; A="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.1658013333333333; H=0.497404; I=0.09737862910630474; J=0.09737862910630474; M=160.56988;          N=481.70964;  O=0.297245424023111;   P=0.297245424023111;   Q=26.62270019717333;  R=239.60430177456;   S=0.02894535189949281;  T=0.02894535189949281  }
)

foreach ($item in $rowData) {
  $ws.Cells.Item($item.row, 1).Value  = $item.A
  $ws.Cells.Item($item.row, 5).Value  = $item.E
  $ws.Cells.Item($item.row, 6).Value  = $item.F
  $ws.Cells.Item($item.row, 7).Value  = $item.G
  $ws.Cells.Item($item.row, 8).Value  = $item.H
  $ws.Cells.Item($item.row, 9).Value  = $item.I
  $ws.Cells.Item($item.row, 10).Value = $item.J
  $ws.Cells.Item($item.row, 13).Value = $item.M
  $ws.Cells.Item($item.row, 14).Value = $item.N
  $ws.Cells.Item($item.row, 15).Value = $item.O
  $ws.Cells.Item($item.row, 16).Value = $item.P
  $ws.Cells.Item($item.row, 17).Value = $item.Q
  $ws.Cells.Item($item.row, 18).Value = $item.R
  $ws.Cells.Item($item.row, 19).Value = $item.S
  $ws.Cells.Item($item.row, 20).Value = $item.T
}

Write-Output "Updated rows 2-13 with new TPM-derived values"
